$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 214 (pushes L16/anc.. and everything below down by one)
$ws.Rows.Item(214).Insert()

# New "Lymphopenia" variable row (group WNL/HI)
$ws.Range("A214").Value2 = "L15a"
$ws.Range("B214").Value2 = "lymphopenia"
$ws.Range("C214").Value2 = "Laboratory"
$ws.Range("D214").Value2 = "Lymphopenia"
$ws.Range("E214").Value2 = "Not lymphopenic; Lymphopenic; Not drawn/Not available; Unknown"

# Description column wraps to two lines for this row
$ws.Rows.Item(214).RowHeight = 31

# Grow the table / autofilter / dimension to include the new row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E299"))

# Reflect where the author was looking/selecting when done
$ws.Range("E215").Select()
